$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (2..5) used to hold species names ("Teh oolong", "Teh hitam",
# "Teh hijau ", "Teh putih"); they're replaced with plain numeric codes.
$ws.Range("A2").Value = 3
$ws.Range("A3").Value = 4
$ws.Range("A4").Value = 8
$ws.Range("A5").Value = 9

# B1 already reads "Target " and B2:B5 already hold 10/15/45/30 - those
# values are unchanged. What changes is formatting: the whole A1:B5 table
# now shares the centered alignment that used to be applied to B2:B5 only.
$ws.Range("A1:B5").HorizontalAlignment = -4108

# Move the active selection to C3.
$ws.Range("C3").Select()
